$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting existing rows 172:228 down to 173:229
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with its data
$ws.Range("A172").Value = 5
$ws.Range("B172").Value = "Macroferia Regional de Talca"
$ws.Range("C172").Value = "Maule"
$ws.Range("D172").Value = 44524
$ws.Range("D172").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E172").Value = 7
$ws.Range("F172").Value = 100114013
$ws.Range("G172").Value = "Zanahoria"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 400
$ws.Range("K172").Value = 9000
$ws.Range("L172").Value = 9000
$ws.Range("M172").Value = 9000
$ws.Range("N172").Value = '$/saco 20 kilos'
$ws.Range("O172").Value = "Provincia del Elquí"
$ws.Range("P172").Value = 450
$ws.Range("Q172").Value = 20
$ws.Range("R172").Value = "Hortaliza"
